# OpenTBS demo workbook update:
#  - Rework the "notes" section text and add a new merge-block example
#    (header row + data row) below it.
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Title cell (B2) now reuses the text that used to be the section
#        title further down; its existing bold Times-New-Roman style is
#        left untouched. -------------------------------------------------
$ws.Range("B2").Value = "OpenTBS demo"

# --- 2. The "notes" list (B12:B16) gets new wording and a new font
#        colour (bold, dark red-ish accent colour). -----------------------
$notes = @(
    "You may consider the following before building your own Microsoft Excel template:",
    "Merging Microsoft Excel templates with OpenTBS has several limitations because of the OpenXML format for Excel.",
    "* Formulas won't work because OpenTBS needs to convert cell positions from aboslute to relative in order to have a constistent merged sheet.",
    "* Formulas may also make troubles because they are saved twice in the sheet:  one for the expression, and one for the instant result.",
    "* Changing picture (using ope=changepic)  because drawing information are saved in another XML sub-file."
)
for ($i = 0; $i -lt $notes.Length; $i++) {
    $row = 12 + $i
    $cell = $ws.Range("B$row")
    $cell.Value = $notes[$i]
    $cell.Font.Bold = $true
    $cell.Font.ThemeColor = 6
}

# Drop the old B17 content (no longer used) since the notes list is now
# only 5 rows (12-16) instead of 7 (12-18).
$ws.Range("B17").ClearContents()

# --- 3. New "Example #1" section title (B18), same style previously
#        used for the "OpenTBS demo" title (bold 16pt Times New Roman). --
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B18").Value = "Example #1: merging data with rows"

# --- 4. Header row (B20:D20) + data/template row (B21:D21) for the new
#        merge-block example table. ---------------------------------------
$ws.Range("B20").Value = "First Name"
$ws.Range("C20").Value = "Name"
$ws.Range("D20").Value = "Membership number"

$ws.Range("B21").Value = "[a.firstname;block=row]"
$ws.Range("C21").Value = "[a.name]"
$ws.Range("D21").Value = "[a.number]"

$headerRange = $ws.Range("B20:D20")
$headerRange.Interior.Color = 14211288
$headerRange.Interior.Pattern = 1
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$dataRange = $ws.Range("B21:D21")
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# --- 5. Column widths for the new table columns. --------------------------
$ws.Columns.Item(2).ColumnWidth = 15.86
$ws.Columns.Item(3).ColumnWidth = 12.1875
$ws.Columns.Item(4).ColumnWidth = 19

# --- 6. Selection / active cell tidy-up (cosmetic, mirrors template). -----
$ws.Range("B16").Select() | Out-Null
